# Reorders the "Recorded By" (column G) contributor list on the
# "Session Analysis Results" sheet so that "System" (or, when "System"
# is not present, "admin@admin.com") is listed first among the
# comma-separated recorder names, matching the upstream normalization.
# Comparisons are case-sensitive (exact token "System", not "system").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }
    if ($raw -isnot [string]) { continue }
    if ($raw -notmatch ",") { continue }

    $parts = @($raw -split "," | ForEach-Object { $_.Trim() })

    $hasSystem = $false
    $hasAdmin = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
        if ($p.Equals("admin@admin.com")) { $hasAdmin = $true }
    }

    $newOrder = $null

    if ($hasSystem) {
        if (-not $parts[0].Equals("System")) {
            $rest = @($parts | Where-Object { -not $_.Equals("System") })
            $newOrder = @("System") + $rest
        }
    } elseif ($hasAdmin) {
        if (-not $parts[0].Equals("admin@admin.com")) {
            $rest = @($parts | Where-Object { -not $_.Equals("admin@admin.com") })
            $newOrder = @("admin@admin.com") + $rest
        }
    }

    if ($newOrder -ne $null) {
        $cell.Value2 = [string]::Join(", ", $newOrder)
    }
}
